$d = $word.ActiveDocument

# The document paragraph holds a header line followed by several
# "<br/>N" runs (1, 2, 3, 4, 5). The edit removes the trailing "4" and
# "5" entries (and their preceding line breaks), leaving the paragraph
# ending right after "3".

$p = $d.Paragraphs.Item(1)

$finder = $p.Range.Duplicate
$found = $finder.Find.Execute("4", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Include the line break that precedes "4" and extend through the
    # end of the paragraph (covering the break + "5" that follow too),
    # but keep the paragraph mark itself intact.
    $deleteStart = $finder.Start - 1
    $deleteEnd = $p.Range.End - 1

    $deleteRange = $d.Range($deleteStart, $deleteEnd)
    $deleteRange.Delete()
}
